$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.044256929020696
$ws.Range("D2").Value = 1.060621967228819
$ws.Range("E2").Value = 1.053262906865078
$ws.Range("F2").Value = 1.066819916247952
$ws.Range("I2").Value = 1.046958685030606
$ws.Range("J2").Value = 1.04932298524398
$ws.Range("K2").Value = 1.063348246543967
$ws.Range("L2").Value = 1.056009359330313
$ws.Range("M2").Value = 1.069529438952142
$ws.Range("N2").Value = 1.020213706455956

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.045692502453099
$ws.Range("D3").Value = 1.061468951822094
$ws.Range("E3").Value = 1.054470043189258
$ws.Range("F3").Value = 1.067951751476877
$ws.Range("I3").Value = 1.04722131680939
$ws.Range("J3").Value = 1.050403334530941
$ws.Range("K3").Value = 1.064009610429053
$ws.Range("L3").Value = 1.057028495864109
$ws.Range("M3").Value = 1.070476152257948
$ws.Range("N3").Value = 1.020579638937778

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.046620592230718
$ws.Range("D4").Value = 1.062015011312353
$ws.Range("E4").Value = 1.055250364797497
$ws.Range("F4").Value = 1.068682560804957
$ws.Range("I4").Value = 1.047388672117753
$ws.Range("J4").Value = 1.051101155815873
$ws.Range("K4").Value = 1.064434919291151
$ws.Range("L4").Value = 1.057686602630173
$ws.Range("M4").Value = 1.071086590525695
$ws.Range("N4").Value = 1.020815802852268

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.047010569526847
$ws.Range("D5").Value = 1.062244097748087
$ws.Range("E5").Value = 1.055578229362672
$ws.Range("F5").Value = 1.068989421079715
$ws.Range("I5").Value = 1.047458409799171
$ws.Range("J5").Value = 1.051394227240265
$ws.Range("K5").Value = 1.064613088593896
$ws.Range("L5").Value = 1.057962952193285
$ws.Range("M5").Value = 1.071342706570365
$ws.Range("N5").Value = 1.020914939124162

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.047076037343073
$ws.Range("D6").Value = 1.062282534424534
$ws.Range("E6").Value = 1.055633268648592
$ws.Range("F6").Value = 1.069040922526737
$ws.Range("I6").Value = 1.047470082804403
$ws.Range("J6").Value = 1.051443418161944
$ws.Range("K6").Value = 1.06464296705538
$ws.Range("L6").Value = 1.058009333907245
$ws.Range("M6").Value = 1.071385679607537
$ws.Range("N6").Value = 1.020931575950494

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.046625803873057
$ws.Range("D7").Value = 1.062018074249664
$ws.Range("E7").Value = 1.055254746453572
$ws.Range("F7").Value = 1.068686662545649
$ws.Range("I7").Value = 1.047389606385915
$ws.Range("J7").Value = 1.051105072996201
$ws.Range("K7").Value = 1.064437302474614
$ws.Range("L7").Value = 1.057690296474551
$ws.Range("M7").Value = 1.07109001477083
$ws.Range("N7").Value = 1.020817128092364

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.044742261167431
$ws.Range("D8").Value = 1.060908623522255
$ws.Range("E8").Value = 1.053671026920153
$ws.Range("F8").Value = 1.067202750502182
$ws.Range("I8").Value = 1.047047978142473
$ws.Range("J8").Value = 1.049688352565702
$ws.Range("K8").Value = 1.063572304322119
$ws.Range("L8").Value = 1.05635406108331
$ws.Range("M8").Value = 1.069849831330347
$ws.Range("N8").Value = 1.020337503869923

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.041416687206854
$ws.Range("D9").Value = 1.05893829779798
$ws.Range("E9").Value = 1.050874219850268
$ws.Range("F9").Value = 1.064575821889258
$ws.Range("I9").Value = 1.046426169455007
$ws.Range("J9").Value = 1.047182279447396
$ws.Range("K9").Value = 1.062027819804771
$ws.Range("L9").Value = 1.053989027430097
$ws.Range("M9").Value = 1.067647915362229
$ws.Range("N9").Value = 1.019487554298124

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.039194917290081
$ws.Range("D10").Value = 1.057614372624081
$ws.Range("E10").Value = 1.049005379908639
$ws.Range("F10").Value = 1.062816251956484
$ws.Range("I10").Value = 1.045998278667867
$ws.Range("J10").Value = 1.045504866172075
$ws.Range("K10").Value = 1.060984482387857
$ws.Range("L10").Value = 1.052405146378307
$ws.Range("M10").Value = 1.066168697942456
$ws.Range("N10").Value = 1.018917628625148

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.038231672628865
$ws.Range("D11").Value = 1.057038621166413
$ws.Range("E11").Value = 1.048195084161677
$ws.Range("F11").Value = 1.062052337118542
$ws.Range("I11").Value = 1.045809823045229
$ws.Range("J11").Value = 1.044776888625628
$ws.Range("K11").Value = 1.060529444111991
$ws.Range("L11").Value = 1.051717559557442
$ws.Range("M11").Value = 1.065525473222629
$ws.Range("N11").Value = 1.018670047242491

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.037873693363967
$ws.Range("D12").Value = 1.056824386777333
$ws.Range("E12").Value = 1.047893938455852
$ws.Range("F12").Value = 1.061768279911262
$ws.Range("I12").Value = 1.045739344310483
$ws.Range("J12").Value = 1.044506233529811
$ws.Range("K12").Value = 1.060359930303192
$ws.Range("L12").Value = 1.051461891576002
$ws.Range("M12").Value = 1.065286140169034
$ws.Range("N12").Value = 1.018577962885475

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.037950489712121
$ws.Range("D13").Value = 1.056870357761456
$ws.Range("E13").Value = 1.04795854280164
$ws.Range("F13").Value = 1.061829225012376
$ws.Range("I13").Value = 1.045754483877833
$ws.Range("J13").Value = 1.044564301433936
$ws.Range("K13").Value = 1.060396313895962
$ws.Range("L13").Value = 1.05151674539809
$ws.Range("M13").Value = 1.065337496558073
$ws.Range("N13").Value = 1.018597720819005

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.038202085809848
$ws.Range("D14").Value = 1.057020920137404
$ws.Range("E14").Value = 1.048170194769976
$ws.Range("F14").Value = 1.062028863107052
$ws.Range("I14").Value = 1.045804007005018
$ws.Range("J14").Value = 1.044754521353942
$ws.Range("K14").Value = 1.060515442100572
$ws.Range("L14").Value = 1.051696431443977
$ws.Range("M14").Value = 1.065505698274706
$ws.Range("N14").Value = 1.018662438009988

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.038357077454492
$ws.Range("D15").Value = 1.057113636891346
$ws.Range("E15").Value = 1.048300578439172
$ws.Range("F15").Value = 1.062151826180709
$ws.Range("I15").Value = 1.045834456485388
$ws.Range("J15").Value = 1.044871688622106
$ws.Range("K15").Value = 1.060588775631455
$ws.Range("L15").Value = 1.051807106369841
$ws.Range("M15").Value = 1.065609278383404
$ws.Range("N15").Value = 1.018702296249588

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.039258818754984
$ws.Range("D16").Value = 1.057652530899927
$ws.Range("E16").Value = 1.049059133549337
$ws.Range("F16").Value = 1.0628669078488
$ws.Range("I16").Value = 1.046010718866542
$ws.Range("J16").Value = 1.0455531445749
$ws.Range("K16").Value = 1.061014612826189
$ws.Range("L16").Value = 1.052450741931562
$ws.Range("M16").Value = 1.06621132921387
$ws.Range("N16").Value = 1.018934042811976

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.039824130444985
$ws.Range("D17").Value = 1.057989898954742
$ws.Range("E17").Value = 1.049534664195101
$ws.Range("F17").Value = 1.063314919367031
$ws.Range("I17").Value = 1.04612043274758
$ws.Range("J17").Value = 1.045980159937933
$ws.Range("K17").Value = 1.061280853645817
$ws.Range("L17").Value = 1.052854004731401
$ws.Range("M17").Value = 1.066588251100794
$ws.Range("N17").Value = 1.019079196159881

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.040153751497873
$ws.Range("D18").Value = 1.058186440560097
$ws.Range("E18").Value = 1.049811929669098
$ws.Range("F18").Value = 1.063576043252059
$ws.Range("I18").Value = 1.046184120604764
$ws.Range("J18").Value = 1.046229072211916
$ws.Range("K18").Value = 1.061435832374078
$ws.Range("L18").Value = 1.053089051944642
$ws.Range("M18").Value = 1.066807841625062
$ws.Range("N18").Value = 1.01916378458095

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.040266124330338
$ws.Range("D19").Value = 1.058253415596879
$ws.Range("E19").Value = 1.049906452599256
$ws.Range("F19").Value = 1.063665047001165
$ws.Range("I19").Value = 1.046205784570739
$ws.Range("J19").Value = 1.046313918103551
$ws.Range("K19").Value = 1.061488622727201
$ws.Range("L19").Value = 1.053169168393166
$ws.Range("M19").Value = 1.066882672031806
$ws.Range("N19").Value = 1.019192614017435

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.039763489864275
$ws.Range("D20").Value = 1.057953727336821
$ws.Range("E20").Value = 1.049483654993683
$ws.Range("F20").Value = 1.063266872025804
$ws.Range("I20").Value = 1.04610869318368
$ws.Range("J20").Value = 1.04593436166394
$ws.Range("K20").Value = 1.061252321123082
$ws.Range("L20").Value = 1.052810755973726
$ws.Range("M20").Value = 1.066547838033725
$ws.Range("N20").Value = 1.019063630561701

# Row 21
$ws.Range("B21").Value = 1.019999999999999
$ws.Range("C21").Value = 1.038128002252843
$ws.Range("D21").Value = 1.056976593625945
$ws.Range("E21").Value = 1.048107873150825
$ws.Range("F21").Value = 1.061970083116618
$ws.Range("I21").Value = 1.045789436877654
$ws.Range("J21").Value = 1.044698513316044
$ws.Range("K21").Value = 1.060480375417422
$ws.Range("L21").Value = 1.051643525800554
$ws.Range("M21").Value = 1.065456178424847
$ws.Range("N21").Value = 1.018643383778925

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.037098619323264
$ws.Range("D22").Value = 1.056360062373695
$ws.Range("E22").Value = 1.047241903213353
$ws.Range("F22").Value = 1.061152972850985
$ws.Range("I22").Value = 1.045585942245676
$ws.Range("J22").Value = 1.043920027213925
$ws.Range("K22").Value = 1.059992173202344
$ws.Range("L22").Value = 1.050908091568485
$ws.Range("M22").Value = 1.064767430675399
$ws.Range("N22").Value = 1.018378453636014

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.037644419895684
$ws.Range("D23").Value = 1.05668710335031
$ws.Range("E23").Value = 1.047701062783591
$ws.Range("F23").Value = 1.061586307134029
$ws.Range("I23").Value = 1.045694080997647
$ws.Range("J23").Value = 1.044332857235625
$ws.Range("K23").Value = 1.060251249062898
$ws.Range("L23").Value = 1.05129810740477
$ws.Range("M23").Value = 1.065132775321804
$ws.Range("N23").Value = 1.018518965363653

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.039790891090301
$ws.Range("D24").Value = 1.05797007247438
$ws.Range("E24").Value = 1.049506704176376
$ws.Range("F24").Value = 1.063288583149974
$ws.Range("I24").Value = 1.046113998734226
$ws.Range("J24").Value = 1.045955056422879
$ws.Range("K24").Value = 1.061265214715757
$ws.Range("L24").Value = 1.05283029874873
$ws.Range("M24").Value = 1.066566099768017
$ws.Range("N24").Value = 1.019070664223777

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.042277235533223
$ws.Range("D25").Value = 1.059449497804003
$ws.Range("E25").Value = 1.051598004469966
$ws.Range("F25").Value = 1.065256394283422
$ws.Range("I25").Value = 1.046589273048007
$ws.Range("J25").Value = 1.047831323514376
$ws.Range("K25").Value = 1.06242951232218
$ws.Range("L25").Value = 1.054601698990248
$ws.Range("M25").Value = 1.068219139779789
$ws.Range("N25").Value = 1.019707861680131
